# Update Excel file after daily scrape

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths. The COM ColumnWidth setter/XML-serialized width are
# offset by 5/6 of a character relative to each other in this runtime, so we
# subtract that fixed offset here to land exactly on the target XML widths
# (67, 81, 17, 42) after round-tripping through the COM property.
$widthOffset = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 67 - $widthOffset
$ws.Columns.Item(4).ColumnWidth = 81 - $widthOffset
$ws.Columns.Item(6).ColumnWidth = 17 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 42 - $widthOffset

# New data for rows 2-18 (A..H)
$data = @(
    @("1327832", "https://aiesec.org/opportunity/global-talent/1327832", "Accelerate Serbia | Intern – Technical Support for eBar Software", "Београд, Србија", "No", "0 applicants", "9 - 12 Weeks", "eBar Software doo"),
    @("1327831", "https://aiesec.org/opportunity/global-talent/1327831", "[Accelerate Serbia] IT Consultant", "Belgrade, Serbia", "No", "0 applicants", "9 - 12 Weeks", "ITSM Solutions and Integrations"),
    @("1327825", "https://aiesec.org/opportunity/global-talent/1327825", "Project Management Trainee", "Panamá, Provincia de Panamá, Panamá", "No", "3 applicants", "6 - 18 Months", "HILTI Panama"),
    @("1327824", "https://aiesec.org/opportunity/global-talent/1327824", "[IMPACT FORTALEZA] Business Development", "Castanhal, Pará, Brasil", "No", "1 applicant", "6 - 18 Months", "Petruz Fruity"),
    @("1327823", "https://aiesec.org/opportunity/global-talent/1327823", "[IMPACT FORTALEZA] Chemical Engiineering", "Castanhal, Pará, Brasil", "No", "0 applicants", "6 - 18 Months", "Petruz Fruity"),
    @("1327819", "https://aiesec.org/opportunity/global-talent/1327819", "[EUROPE ONLY] HR Opportunities with Dutch", "Bucharest, Romania", "No", "2 applicants", "6 - 18 Months", "Accenture Romania"),
    @("1327818", "https://aiesec.org/opportunity/global-talent/1327818", "[EUROPE ONLY] Recruiting Analyst with Polish", "Bucharest, Romania", "No", "0 applicants", "6 - 18 Months", "Accenture Romania"),
    @("1327817", "https://aiesec.org/opportunity/global-talent/1327817", "[EUROPE ONLY] HR Admin with Polish", "Bucharest, Romania", "No", "0 applicants", "6 - 18 Months", "Accenture Romania"),
    @("1327815", "https://aiesec.org/opportunity/global-talent/1327815", "Computer and AI Coordinator", "London, UK", "No", "24 applicants", "6 - 18 Months", "Capital Care Homes"),
    @("1327813", "https://aiesec.org/opportunity/global-talent/1327813", "Nursery Spanish Practitioner", "Ashby-de-la-Zouch LE65, UK", "No", "1 applicant", "6 - 18 Months", "Bilingual Day Nursery and Preschool Ltd"),
    @("1327811", "https://aiesec.org/opportunity/global-talent/1327811", "Software Engineering Intern", "Colombo, Sri Lanka", "No", "4 applicants", "3 - 6 Months", "Envision Circle (Pvt) Ltd"),
    @("1327778", "https://aiesec.org/opportunity/global-talent/1327778", "Digital Content & Stakeholder Engagement Intern", "Colombo, Sri Lanka", "No", "3 applicants", "6 - 18 Months", "Solutions Ground (Pvt) Ltd"),
    @("1327658", "https://aiesec.org/opportunity/global-talent/1327658", "DevOps Engineer", "El-Kom El-Ahmar, Shibin el-Qanater, Al-Qalyubia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Etolv"),
    @("1327399", "https://aiesec.org/opportunity/global-talent/1327399", "Architectural Designer / 3D Visualizer", "Shebeen El-Kom, Qism Shebeen El-Kom, Shibin el Kom, Menofia Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "Building Consultant Office"),
    @("1327397", "https://aiesec.org/opportunity/global-talent/1327397", "Site Execution Engineer", "Shebeen El-Kom, Qism Shebeen El-Kom, Shibin el Kom, Menofia Governorate, Egypt", "No", "1 applicant", "9 - 12 Weeks", "Building Consultant Office"),
    @("1326701", "https://aiesec.org/opportunity/global-talent/1326701", "Software Developer", "Berlin, Germany", "No", "210 applicants", "6 - 18 Months", "code4business Shareholder GmbH"),
    @("1307425", "https://aiesec.org/opportunity/global-talent/1307425", "Service executive II", "Naucalpan de Juárez, Mexico", "No", "22 applicants", "6 - 18 Months", "Segmenta S.C.")
)

$rowIndex = 2
foreach ($row in $data) {
    # Column A holds numeric-looking opportunity IDs that must stay text,
    # matching the original inline-string cell type.
    $ws.Cells.Item($rowIndex, 1).NumberFormat = "@"
    for ($col = 1; $col -le 8; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
